$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 401.75
$ws.Range("J96").Value = 220
$ws.Range("L96").Value = 660
$ws.Range("N96").Value = -3406

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3600.426
$ws.Range("I132").Value = 3819.439
$ws.Range("J132").Value = 2909.6924
$ws.Range("K132").Value = 11458.317
$ws.Range("L132").Value = 8729.0772
$ws.Range("M132").Value = -8928.316999999999
$ws.Range("N132").Value = -13789.0772

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3333.1482
$ws.Range("I20").Value = 3745.8667
$ws.Range("J20").Value = 2817.25
$ws.Range("K20").Value = 3745.8667
$ws.Range("L20").Value = 2817.25
$ws.Range("M20").Value = -3498.8667
$ws.Range("N20").Value = -3311.25
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15335.546
$ws.Range("J50").Value = 16660.8
$ws.Range("L50").Value = 16660.8
$ws.Range("N50").Value = -17910.8
$ws.Range("H54").Value = 49998.332
$ws.Range("J54").Value = 49998.332
$ws.Range("L54").Value = 49998.332
$ws.Range("N54").Value = -51314.332
$ws.Range("H58").Value = 9524728
$ws.Range("I58").Value = 987.6539
$ws.Range("J58").Value = 37037756
$ws.Range("K58").Value = 987.6539
$ws.Range("L58").Value = 37037756
$ws.Range("M58").Value = -784.6539
$ws.Range("N58").Value = -37038162
$ws.Range("H69").Value = 5960.6665
$ws.Range("I69").Value = 5960.6665
$ws.Range("K69").Value = 5960.6665
$ws.Range("M69").Value = -5211.6665
$ws.Range("H72").Value = 5960.6665
$ws.Range("I72").Value = 5960.6665
$ws.Range("K72").Value = 17881.9995
$ws.Range("M72").Value = -14137.9995
$ws.Range("H99").Value = 590977.8
$ws.Range("I99").Value = 1669035.1
$ws.Range("J99").Value = 2946.5454
$ws.Range("K99").Value = 1669035.1
$ws.Range("L99").Value = 2946.5454
$ws.Range("M99").Value = -1667537.1
$ws.Range("N99").Value = -5942.5454
$ws.Range("H126").Value = 590977.8
$ws.Range("I126").Value = 1669035.1
$ws.Range("J126").Value = 2946.5454
$ws.Range("K126").Value = 5007105.300000001
$ws.Range("L126").Value = 8839.636200000001
$ws.Range("M126").Value = -5004635.300000001
$ws.Range("N126").Value = -13779.6362
$ws.Range("H132").Value = 1978.8206
$ws.Range("I132").Value = 1913.5294
$ws.Range("J132").Value = 2422.8
$ws.Range("K132").Value = 5740.5882
$ws.Range("L132").Value = 7268.400000000001
$ws.Range("M132").Value = -3210.5882
$ws.Range("N132").Value = -12328.4
$ws.Range("H134").Value = 1075.0435
$ws.Range("I134").Value = 1065.3024
$ws.Range("J134").Value = 1214.6666
$ws.Range("K134").Value = 3195.9072
$ws.Range("L134").Value = 3643.9998
$ws.Range("M134").Value = -660.9072000000001
$ws.Range("N134").Value = -8713.9998
$ws.Range("H136").Value = 9524728
$ws.Range("I136").Value = 987.6539
$ws.Range("J136").Value = 37037756
$ws.Range("K136").Value = 2962.9617
$ws.Range("L136").Value = 111113268
$ws.Range("M136").Value = -412.9616999999998
$ws.Range("N136").Value = -111118368

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 8185.7144
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 9460
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 9460
$ws.Range("M53").Value = -4369
$ws.Range("N53").Value = -10722
$ws.Range("H80").Value = 4978
$ws.Range("I80").Value = 3453.3333
$ws.Range("J80").Value = 6121.5
$ws.Range("K80").Value = 3453.3333
$ws.Range("L80").Value = 6121.5
$ws.Range("M80").Value = -2455.3333
$ws.Range("N80").Value = -8117.5
$ws.Range("H83").Value = 4978
$ws.Range("I83").Value = 3453.3333
$ws.Range("J83").Value = 6121.5
$ws.Range("K83").Value = 17266.6665
$ws.Range("L83").Value = 30607.5
$ws.Range("M83").Value = -12274.6665
$ws.Range("N83").Value = -40591.5
$ws.Range("H113").Value = 27779906
$ws.Range("I113").Value = 125005000
$ws.Range("J113").Value = 1307.1428
$ws.Range("K113").Value = 125005000
$ws.Range("L113").Value = 1307.1428
$ws.Range("M113").Value = -125002830
$ws.Range("N113").Value = -5647.1428
$ws.Range("H132").Value = 27441.281
$ws.Range("I132").Value = 34900.266
$ws.Range("J132").Value = 2578
$ws.Range("K132").Value = 104700.798
$ws.Range("L132").Value = 7734
$ws.Range("M132").Value = -102170.798
$ws.Range("N132").Value = -12794

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2476.6667
$ws.Range("I22").Value = 3220
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 3220
$ws.Range("L22").Value = 990
$ws.Range("M22").Value = -2925
$ws.Range("N22").Value = -1580
$ws.Range("H27").Value = 2476.6667
$ws.Range("I27").Value = 3220
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 3220
$ws.Range("L27").Value = 990
$ws.Range("M27").Value = -3113
$ws.Range("N27").Value = -1204
$ws.Range("H93").Value = 2704593.2
$ws.Range("I93").Value = 3862847.2
$ws.Range("J93").Value = 2001
$ws.Range("K93").Value = 3862847.2
$ws.Range("L93").Value = 2001
$ws.Range("M93").Value = -3861599.2
$ws.Range("N93").Value = -4497
$ws.Range("H132").Value = 7256.5757
$ws.Range("I132").Value = 12269.294
$ws.Range("J132").Value = 1930.5625
$ws.Range("K132").Value = 36807.882
$ws.Range("L132").Value = 5791.6875
$ws.Range("M132").Value = -34277.882
$ws.Range("N132").Value = -10851.6875
$ws.Range("H136").Value = 4518.0625
$ws.Range("I136").Value = 5298.826
$ws.Range("J136").Value = 2522.7778
$ws.Range("K136").Value = 15896.478
$ws.Range("L136").Value = 7568.3334
$ws.Range("M136").Value = -13346.478
$ws.Range("N136").Value = -12668.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 317750
$ws.Range("J54").Value = 6000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7040
$ws.Range("H81").Value = 2054.4546
$ws.Range("I81").Value = 2079.9
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 4159.8
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = -3098.8
$ws.Range("N81").Value = -5722
$ws.Range("H84").Value = 2054.4546
$ws.Range("I84").Value = 2079.9
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 20799
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -15495
$ws.Range("N84").Value = -28608
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -259
$ws.Range("H132").Value = 1340.5927
$ws.Range("I132").Value = 1247.84
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3743.52
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1213.52
$ws.Range("N132").Value = -12560
